# Applies the cryptos.xlsx update described by the commit diff.
# Column D values are numeric-looking text (e.g. '70.683.10', '0.102') that must
# stay as text, so they're written with a leading apostrophe (force-text, like
# typing '70.683.10 directly into Excel) to avoid Excel's automatic number coercion.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 4).Value = "'70.683.10"
$ws.Cells.Item(2, 5).Value = "  +1.83%  "

$ws.Cells.Item(3, 4).Value = "'3.577.45"
$ws.Cells.Item(3, 5).Value = "  +0.98%  "

$ws.Cells.Item(4, 5).Value = "  -0.06%  "

$ws.Cells.Item(5, 4).Value = "'607.45"
$ws.Cells.Item(5, 5).Value = "  +4.50%  "

$ws.Cells.Item(6, 4).Value = "'174.29"
$ws.Cells.Item(6, 5).Value = "  +0.98%  "

$ws.Cells.Item(7, 4).Value = "'3.570.95"
$ws.Cells.Item(7, 5).Value = "  +0.96%  "

$ws.Cells.Item(8, 5).Value = "  +0.51%  "

$ws.Cells.Item(9, 5).Value = "  -0.02%  "

$ws.Cells.Item(10, 5).Value = "  +3.38%  "

$ws.Cells.Item(11, 4).Value = "'7.42"
$ws.Cells.Item(11, 5).Value = "  +9.79%  "

$ws.Cells.Item(12, 5).Value = "  +0.32%  "

$ws.Cells.Item(13, 4).Value = "'46.98"
$ws.Cells.Item(13, 5).Value = "  -1.13%  "

$ws.Cells.Item(14, 5).Value = "  +0.89%  "

$ws.Cells.Item(15, 4).Value = "'4.156.97"
$ws.Cells.Item(15, 5).Value = "  +1.00%  "

$ws.Cells.Item(16, 4).Value = "'8.44"
$ws.Cells.Item(16, 5).Value = "  -1.60%  "

$ws.Cells.Item(17, 4).Value = "'618.40"
$ws.Cells.Item(17, 5).Value = "  -2.06%  "

$ws.Cells.Item(18, 4).Value = "'3.601.18"
$ws.Cells.Item(18, 5).Value = "  +1.35%  "

$ws.Cells.Item(19, 4).Value = "'70.814.51"
$ws.Cells.Item(19, 5).Value = "  +1.94%  "

$ws.Cells.Item(20, 5).Value = "  -2.37%  "

$ws.Cells.Item(21, 5).Value = "  +0.12%  "

$ws.Cells.Item(22, 5).Value = "  -0.56%  "

$ws.Cells.Item(23, 5).Value = "  -16.59%  "

$ws.Cells.Item(24, 4).Value = "'16.17"
$ws.Cells.Item(24, 5).Value = "  +0.75%  "

$ws.Cells.Item(25, 4).Value = "'97.74"
$ws.Cells.Item(25, 5).Value = "  -0.13%  "

$ws.Cells.Item(26, 5).Value = "  -0.19%  "

$ws.Cells.Item(27, 5).Value = "  +0.01%  "

$ws.Cells.Item(28, 5).Value = "  +0.23%  "

$ws.Cells.Item(29, 2).Value = "RenderToken"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(29, 4).Value = "'9.28"
$ws.Cells.Item(29, 5).Value = "  -0.72%  "

$ws.Cells.Item(30, 2).Value = "EthereumClassic"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(30, 4).Value = "'33.57"
$ws.Cells.Item(30, 5).Value = "  +1.80%  "

$ws.Cells.Item(31, 4).Value = "'8.51"
$ws.Cells.Item(31, 5).Value = "  -0.98%  "

$ws.Cells.Item(32, 4).Value = "'3.06"
$ws.Cells.Item(32, 5).Value = "  -3.43%  "

$ws.Cells.Item(33, 4).Value = "'7.07"
$ws.Cells.Item(33, 5).Value = "  +0.63%  "

$ws.Cells.Item(34, 5).Value = "  -2.93%  "

$ws.Cells.Item(35, 2).Value = "dogwifhat"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Cells.Item(35, 4).Value = "'3.83"
$ws.Cells.Item(35, 5).Value = "  +8.89%  "

$ws.Cells.Item(36, 2).Value = "Hedera"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(36, 4).Value = "'0.102"
$ws.Cells.Item(36, 5).Value = "  -1.07%  "

$ws.Cells.Item(37, 2).Value = "Cosmos"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(37, 4).Value = "'10.87"
$ws.Cells.Item(37, 5).Value = "  +0.47%  "

$ws.Cells.Item(38, 2).Value = "VeChain"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(38, 4).Value = "'0.0489"
$ws.Cells.Item(38, 5).Value = "  +6.84%  "

$ws.Cells.Item(39, 2).Value = "OKB"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(39, 4).Value = "'57.34"
$ws.Cells.Item(39, 5).Value = "  -0.07%  "

$ws.Cells.Item(40, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Cells.Item(40, 4).Value = "'1.00"
$ws.Cells.Item(40, 5).Value = "  +0.03%  "

$ws.Cells.Item(41, 2).Value = "Kaspa"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(41, 4).Value = "'0.141"
$ws.Cells.Item(41, 5).Value = "  +3.37%  "

$ws.Cells.Item(42, 2).Value = "Maker"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(42, 4).Value = "'3.393.91"
$ws.Cells.Item(42, 5).Value = "  -0.31%  "

$ws.Cells.Item(43, 2).Value = "TheGraph"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Cells.Item(43, 4).Value = "'0.324"
$ws.Cells.Item(43, 5).Value = "  -1.92%  "

$ws.Cells.Item(44, 2).Value = "ThetaToken"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Cells.Item(44, 4).Value = "'3.00"
$ws.Cells.Item(44, 5).Value = "  +8.67%  "

$ws.Cells.Item(45, 4).Value = "'0.0₃0714"
$ws.Cells.Item(45, 5).Value = "  +1.42%  "

$ws.Cells.Item(46, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(46, 4).Value = "'33.05"
$ws.Cells.Item(46, 5).Value = "  +0.84%  "

$ws.Cells.Item(47, 2).Value = "Fetch.AI"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(47, 4).Value = "'2.67"
$ws.Cells.Item(47, 5).Value = "  +3.78%  "

$ws.Cells.Item(48, 2).Value = "Stellar"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(48, 4).Value = "'0.131"
$ws.Cells.Item(48, 5).Value = "  +0.43%  "

$ws.Cells.Item(49, 2).Value = "Monero"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(49, 4).Value = "'132.85"
$ws.Cells.Item(49, 5).Value = "  -0.29%  "

$ws.Cells.Item(50, 2).Value = "USDe"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Cells.Item(50, 4).Value = "'1.00"
$ws.Cells.Item(50, 5).Value = "  -0.02%  "

$ws.Cells.Item(51, 2).Value = "Cronos"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(51, 4).Value = "'0.148"
$ws.Cells.Item(51, 5).Value = "  +0.62%  "

